$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 12 (task #8): swap/replace task + responsible, and set row height
# ---------------------------------------------------------------------------
$ws.Cells.Item(12, 2).Value = 'Creacion de la estructura static y para crear y vincular el css a las plantillas y creacion del css base'
$ws.Cells.Item(12, 3).Value = 'Sergio'
$ws.Rows.Item(12).RowHeight = 45

# ---------------------------------------------------------------------------
# Row 13 (task #9)
# ---------------------------------------------------------------------------
$ws.Cells.Item(13, 2).Value = 'Modificacion de la pagina principal'
$ws.Cells.Item(13, 3).Value = 'Samuel'

# ---------------------------------------------------------------------------
# Row 14 (task #10)
# ---------------------------------------------------------------------------
$ws.Cells.Item(14, 2).Value = 'Mejoras en la visualizacion del listado haciendolo mediante tablas añadiendo algunos atributos mas para mejorar la comprension, cambios en el estilo de los botones, cambios en el estilo en los links a los detalles, modificado boton de volver al listado y añadido volver a pagina anterior, modificados los botones de volver de las paginas de creacion edicion y eliminacion '

# ---------------------------------------------------------------------------
# Row 15 (task #11), row height 30 -> 60
# ---------------------------------------------------------------------------
$ws.Cells.Item(15, 2).Value = 'Modificacion del formato de los botones de eliminacion y creacion de clases para los que tengan las mismas funciones'
$ws.Cells.Item(15, 3).Value = 'Sergio'
$ws.Rows.Item(15).RowHeight = 60

# ---------------------------------------------------------------------------
# Row 16 (task #12)
# ---------------------------------------------------------------------------
$ws.Cells.Item(16, 2).Value = 'Correcciones leves, actualizacion ReadMe y Organización'

# ---------------------------------------------------------------------------
# Row 17 (task #13)
# ---------------------------------------------------------------------------
$ws.Cells.Item(17, 2).Value = 'Modificacion final al header y adicion de Logo'
$ws.Cells.Item(17, 3).Value = 'Samuel'

# ---------------------------------------------------------------------------
# Row 18 (task #14)
# ---------------------------------------------------------------------------
$ws.Cells.Item(18, 2).Value = 'Responsive y ajustes'

# ---------------------------------------------------------------------------
# Row 19 (task #15), row height unset -> 30
# ---------------------------------------------------------------------------
$ws.Cells.Item(19, 2).Value = 'Modificacion final y mejora  de la pagina principal'
$ws.Cells.Item(19, 3).Value = 'Sergio'
$ws.Rows.Item(19).RowHeight = 30

# ---------------------------------------------------------------------------
# Row 20 (task #16) - previously blank placeholder row, now filled in
# ---------------------------------------------------------------------------
$ws.Cells.Item(20, 1).Value = 16
$ws.Cells.Item(20, 2).Value = 'Entidad Relacion'
$ws.Cells.Item(20, 3).Value = 'Samuel'
$ws.Cells.Item(20, 4).Value = 45753
$ws.Cells.Item(20, 5).Value = 45754
$ws.Cells.Item(20, 6).Value = '✅ Hecho'

# ---------------------------------------------------------------------------
# Row 21 (task #17) - previously blank placeholder row, now filled in
# ---------------------------------------------------------------------------
$ws.Cells.Item(21, 1).Value = 17
$ws.Cells.Item(21, 2).Value = 'Retoques finales para entregar'
$ws.Cells.Item(21, 3).Value = 'Jon'
$ws.Cells.Item(21, 4).Value = 45754
$ws.Cells.Item(21, 5).Value = 45754
$ws.Cells.Item(21, 6).Value = '✅ Hecho'

# ---------------------------------------------------------------------------
# Row 22 - brand new blank placeholder row (styled like row 12, blank values)
# ---------------------------------------------------------------------------
$fmtSrc = $ws.Range("A12:G12")
$row22 = $ws.Range("A22:G22")
$fmtSrc.Copy($row22)
$row22.ClearContents()
$sameStyle = $ws.Range("A22")
$sameStyle.Copy()
$ws.Range("D22:E22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Move the old row 24 (task #6) content down to row 33, removing row 24
# without shifting the rows in between (25-27 must stay put)
# ---------------------------------------------------------------------------
$oldRow = $ws.Range("A24:G24")
$newRow = $ws.Range("A33:G33")
$oldRow.Copy($newRow)
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(24).Insert()
$ws.Rows.Item(33).RowHeight = 45

# ---------------------------------------------------------------------------
# Update the active selection / scroll position
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("I31").Select()
